$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "26.928.65"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.863.53"
Set-TextValue "D4" "1.0000"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue "D5" "305.01"
$ws.Range("E5").Value = "  -0.81%  "
Set-TextValue "D6" "0.9998"
$ws.Range("E6").Value = "  -0.10%  "
Set-TextValue "D7" "0.5063"
$ws.Range("E7").Value = "  +0.30%  "
Set-TextValue "D8" "0.3627"
$ws.Range("E8").Value = "  -3.37%  "
Set-TextValue "D9" "0.07172"
$ws.Range("E9").Value = "  +0.18%  "
Set-TextValue "D10" "0.8970"
$ws.Range("E10").Value = "  +0.76%  "
Set-TextValue "D11" "20.71"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "1.857.28"
$ws.Range("E12").Value = "  -0.92%  "
Set-TextValue "D13" "0.07447"
$ws.Range("E13").Value = "  -1.48%  "
Set-TextValue "D14" "92.62"
$ws.Range("E14").Value = "  +3.63%  "
Set-TextValue "D15" "5.240"
$ws.Range("E16").Value = "  -0.08%  "
Set-TextValue "D17" "0.000008493"
$ws.Range("E17").Value = "  -0.23%  "
Set-TextValue "D18" "14.13"
$ws.Range("E18").Value = "  -0.15%  "
Set-TextValue "D19" "0.9995"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "26.960.19"
$ws.Range("E20").Value = "  -0.75%  "
Set-TextValue "D21" "5.026"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "2.089.27"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  -2.68%  "
Set-TextValue "D24" "6.432"
$ws.Range("E24").Value = "  -0.99%  "
Set-TextValue "D25" "148.10"
$ws.Range("E25").Value = "  -2.09%  "
Set-TextValue "D26" "1.795"
$ws.Range("E26").Value = "  -2.50%  "
Set-TextValue "D27" "17.86"
$ws.Range("E27").Value = "  -0.94%  "
Set-TextValue "D28" "2.064"
$ws.Range("E28").Value = "  -1.54%  "
Set-TextValue "D29" "113.23"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -1.99%  "
Set-TextValue "D31" "4.677"
$ws.Range("E31").Value = "  -0.38%  "
Set-TextValue "D32" "0.09237"
Set-TextValue "D33" "0.05084"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D34" "0.7479"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "2.991"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("E36").Value = "  -0.95%  "
Set-TextValue "D37" "3.272"
$ws.Range("E37").Value = "  +7.52%  "
Set-TextValue "D38" "2.525"
$ws.Range("E38").Value = "  -1.26%  "
Set-TextValue "D39" "0.01998"
$ws.Range("E39").Value = "  -1.94%  "
Set-TextValue "D40" "1.084"
$ws.Range("E40").Value = "  +0.84%  "
Set-TextValue "D41" "0.5422"
$ws.Range("E41").Value = "  +1.00%  "
Set-TextValue "D42" "117.70"
$ws.Range("E42").Value = "  +2.38%  "
Set-TextValue "D43" "6.496"
$ws.Range("E43").Value = "  -1.88%  "
Set-TextValue "D44" "8.574"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  -0.81%  "
Set-TextValue "D46" "0.4664"
$ws.Range("E46").Value = "  +0.16%  "
Set-TextValue "D47" "0.9993"
$ws.Range("E47").Value = "  -0.12%  "
Set-TextValue "D48" "10.09"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("E51").Value = "  -2.50%  "
